$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 75: quarterly "Serie" date 01-04-2021 (as plain text, like the
# other "Serie" entries in column A) plus the 28 numeric data columns.
#
# Column A must end up as a literal text shared-string ("01-04-2021"),
# not an auto-parsed date serial. Excel's normal typed-value parser would
# convert this into a date. To avoid that, force the cell to Text format
# before entering the value, then clear the formatting again afterwards
# so the cell keeps the default (unstyled) look of the rest of column A.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()

$ws.Range("B75").Value = 11383
$ws.Range("C75").Value = -61
$ws.Range("D75").Value = 1456
$ws.Range("E75").Value = 206
$ws.Range("F75").Value = 1250
$ws.Range("G75").Value = 3754
$ws.Range("H75").Value = 3754
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 6095
$ws.Range("K75").Value = 6092
$ws.Range("L75").Value = 4
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = 140
$ws.Range("P75").Value = 118
$ws.Range("Q75").Value = 11265
$ws.Range("R75").Value = 1914
$ws.Range("S75").Value = 1131
$ws.Range("T75").Value = 784
$ws.Range("U75").Value = 9349
$ws.Range("V75").Value = 9350
$ws.Range("W75").Value = -1
$ws.Range("X75").Value = 0
$ws.Range("Y75").Value = 0
$ws.Range("Z75").Value = 0
$ws.Range("AA75").Value = 0
$ws.Range("AB75").Value = 0
$ws.Range("AC75").Value = 3
